$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell originally holds a literal text string (e.g. "304.13" or "1.67%"),
# not a numeric/percentage value. Prefix new values with a leading apostrophe so
# Excel stores them as text (quote-prefixed) instead of auto-converting them to
# numbers/percentages.
$ws.Range('D2').Value = "'304.13"
$ws.Range('E2').Value = "'1.67%"
$ws.Range('D3').Value = "'32.15"
$ws.Range('E3').Value = "'2.42%"
$ws.Range('D4').Value = "'5.015"
$ws.Range('E4').Value = "'-1.34%"
$ws.Range('D5').Value = "'0.07877"
$ws.Range('E5').Value = "'-0.35%"
$ws.Range('D6').Value = "'2.081"
$ws.Range('E6').Value = "'-9.80%"
$ws.Range('D7').Value = "'7.838"
$ws.Range('E7').Value = "'1.26%"
$ws.Range('D8').Value = "'3.837"
$ws.Range('E8').Value = "'-0.79%"
$ws.Range('D9').Value = "'0.9248"
$ws.Range('E9').Value = "'0.80%"
$ws.Range('D10').Value = "'0.1758"
$ws.Range('E10').Value = "'1.75%"
$ws.Range('D11').Value = "'0.07765"
$ws.Range('E11').Value = "'6.24%"
$ws.Range('D12').Value = "'0.08566"
$ws.Range('E12').Value = "'-5.85%"
$ws.Range('D13').Value = "'0.03158"
$ws.Range('E13').Value = "'4.50%"
$ws.Range('D14').Value = "'0.1006"
$ws.Range('E14').Value = "'0.44%"
$ws.Range('D15').Value = "'0.001516"
$ws.Range('E15').Value = "'-0.12%"
$ws.Range('D16').Value = "'0.005662"
$ws.Range('E16').Value = "'-8.10%"
$ws.Range('E17').Value = "'2,108.13%"
$ws.Range('D18').Value = "'3.467"
$ws.Range('E18').Value = "'-0.39%"
$ws.Range('E19').Value = "'-6.51%"
$ws.Range('D20').Value = "'0.3277"
$ws.Range('E20').Value = "'0.16%"
$ws.Range('E21').Value = "'0.48%"
$ws.Range('D22').Value = "'4.282"
$ws.Range('D23').Value = "'0.1858"
$ws.Range('E23').Value = "'9.32%"
$ws.Range('D24').Value = "'0.04580"
$ws.Range('E24').Value = "'-0.89%"
$ws.Range('E25').Value = "'-1.14%"
$ws.Range('D26').Value = "'0.004462"
$ws.Range('E26').Value = "'0.25%"
$ws.Range('D27').Value = "'0.0001247"
$ws.Range('E27').Value = "'3.99%"
$ws.Range('D39').Value = "'0.01740"
$ws.Range('E39').Value = "'0.13%"
$ws.Range('E40').Value = "'4.57%"
$ws.Range('D41').Value = "'0.007462"
$ws.Range('E41').Value = "'7.04%"
$ws.Range('D42').Value = "'0.1364"
$ws.Range('E42').Value = "'0.76%"
$ws.Range('D43').Value = "'0.002354"
$ws.Range('E43').Value = "'7.57%"
$ws.Range('D44').Value = "'0.01041"
$ws.Range('E44').Value = "'9.20%"
$ws.Range('D45').Value = "'0.00006122"
$ws.Range('E45').Value = "'-2.72%"
$ws.Range('E46').Value = "'-0.08%"
$ws.Range('D47').Value = "'0.003095"
$ws.Range('E47').Value = "'-61.20%"
$ws.Range('D48').Value = "'0.8205"
$ws.Range('E48').Value = "'9.81%"
$ws.Range('D49').Value = "'0.00002097"
$ws.Range('E49').Value = "'-0.08%"
$ws.Range('D50').Value = "'0.0001997"
$ws.Range('E50').Value = "'-0.08%"
